$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "00000398"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = 22300001

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "00000399"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = 117640001

# Match styling of A2 (bordered / centered / bold style) for the new A column cells
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
